$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 103695
$ws.Range("D11").Value = 89599.9785

$ws.Range("B12").Value = 107380
$ws.Range("D12").Value = 91551.15549999999
